# rerun with NF1 intronic plp
# Updates the "Neurofibroma plexiform" histology sheet's cpgPLP / no_cpgPLP
# summary statistics to reflect the rerun that reclassifies a case using the
# NF1 intronic PLP variant (cohort split moves from 10/5 to 11/4).

function Set-TextCell($ws, $addr, $val) {
    # Force the cell to be stored as text (matches the source workbook, where
    # these numeric-looking values are shared-string text, not numbers).
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
    $ws.Range($addr).Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Neurofibroma plexiform")

# Number Cohort
$ws.Range("B2").Value = "11 (73.33%)"
$ws.Range("C2").Value = "4 (26.67%)"

# Sex
$ws.Range("B3").Value = "5 (45.5%)"
$ws.Range("C3").Value = "2 (50%)"
$ws.Range("B4").Value = "6 (54.5%)"
$ws.Range("C4").Value = "2 (50%)"

# Race
$ws.Range("B5").Value = "4 (36.4%)"
$ws.Range("C5").Value = "NA (NA%)"
$ws.Range("B6").Value = "2 (18.2%)"
$ws.Range("C6").Value = "1 (25%)"
$ws.Range("B7").Value = "5 (45.5%)"
$ws.Range("C7").Value = "3 (75%)"

# Ethnicity
$ws.Range("B8").Value = "0 (0%)"
$ws.Range("C8").Value = "1 (25%)"
$ws.Range("B9").Value = "10 (90.9%)"
$ws.Range("C9").Value = "3 (75%)"
$ws.Range("B10").Value = "1 (9.1%)"
$ws.Range("C10").Value = "NA (NA%)"

# Ancestry
$ws.Range("B11").Value = "4 (36.4%)"
$ws.Range("C11").Value = "NA (NA%)"
$ws.Range("B12").Value = "2 (18.2%)"
$ws.Range("C12").Value = "1 (25%)"
$ws.Range("B13").Value = "5 (45.5%)"
$ws.Range("C13").Value = "3 (75%)"

# Medians (stored as text in the workbook)
Set-TextCell $ws "B14" "16.95"
Set-TextCell $ws "C14" "12.7"
Set-TextCell $ws "B15" "7.79"
Set-TextCell $ws "C15" "2.7"
Set-TextCell $ws "B16" "7.47"
Set-TextCell $ws "C16" "0.73"
Set-TextCell $ws "B17" "0.25"
Set-TextCell $ws "C17" "0.11"
